$d = $word.ActiveDocument
$t = $d.Tables(1)

# Each table cell holds a single "a÷b=c, d" answer string. Several of the
# old values repeat (e.g. "36÷6=6, 0" occurs twice but needs two different
# replacements), so update by explicit (row, column) cell address rather
# than Find/Replace, which matches every occurrence document-wide.

$t.Cell(1, 1).Range.Text  = "35÷3=11, 2"
$t.Cell(1, 2).Range.Text  = "32÷8=4, 0"
$t.Cell(1, 3).Range.Text  = "45÷4=11, 1"
$t.Cell(1, 4).Range.Text  = "55÷6=9, 1"
$t.Cell(1, 5).Range.Text  = "59÷8=7, 3"

$t.Cell(5, 1).Range.Text  = "92÷3=30, 2"
$t.Cell(5, 2).Range.Text  = "78÷3=26, 0"
$t.Cell(5, 3).Range.Text  = "97÷2=48, 1"
$t.Cell(5, 4).Range.Text  = "19÷7=2, 5"
$t.Cell(5, 5).Range.Text  = "46÷2=23, 0"

$t.Cell(9, 1).Range.Text  = "70÷2=35, 0"
$t.Cell(9, 2).Range.Text  = "32÷3=10, 2"
$t.Cell(9, 3).Range.Text  = "77÷4=19, 1"
$t.Cell(9, 4).Range.Text  = "53÷8=6, 5"
$t.Cell(9, 5).Range.Text  = "96÷7=13, 5"

$t.Cell(13, 1).Range.Text = "35÷2=17, 1"
$t.Cell(13, 2).Range.Text = "41÷9=4, 5"
$t.Cell(13, 3).Range.Text = "88÷3=29, 1"
$t.Cell(13, 4).Range.Text = "74÷7=10, 4"
$t.Cell(13, 5).Range.Text = "91÷2=45, 1"

$t.Cell(17, 1).Range.Text = "20÷8=2, 4"
$t.Cell(17, 2).Range.Text = "17÷5=3, 2"
$t.Cell(17, 3).Range.Text = "33÷9=3, 6"
$t.Cell(17, 4).Range.Text = "32÷8=4, 0"
$t.Cell(17, 5).Range.Text = "45÷4=11, 1"
